$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: update title and link
$ws.Range("D28").Value = "배터리 모델 간단 정리"
$ws.Range("E28").Value = "https://ropiens.tistory.com/253"

# Row 36: update title and link
$ws.Range("D36").Value = "How to Apply AI in Semiconductor Manufacturing: Current Approaches and Case Studies"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/461"

$wb.Save()
